$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, pushing existing rows 17:47 down to 18:48
$ws.Rows("17").Insert()

# Populate the new row 17 with the new weekly price entry
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44915
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112026
$ws.Range("G17").Value = "Haba"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 700
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
